$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the long helper/explanation text from E1 (no longer needed)
$ws.Range("E1").Value = ""

# Rename the AbilitySpecial block header to AbilityValues
$ws.Range("D2").Value = "AbilityValues[{]"

# Replace the zero-padded numeric-looking key placeholders (01..09) with plain 1..9
$ws.Range("E2").Value = "1"
$ws.Range("F2").Value = "2"
$ws.Range("G2").Value = "3"
$ws.Range("H2").Value = "4"
$ws.Range("I2").Value = "5"
$ws.Range("J2").Value = "6"
$ws.Range("K2").Value = "7"
$ws.Range("L2").Value = "8"
$ws.Range("M2").Value = "9"

# Update example value row to include the "damage" key prefix
$ws.Range("E3").Value = "damage 0.1 0.2 0.3 0.4"

# Update the selected cell shown in the sheet view
$ws.Range("D10").Select()
